# Vendors - Monthly.xlsx : add a new vendor ("Atlas Crane Service") to the
# alphabetical list, refresh the "last checked" date in I1, move the
# Envelope "x" marks in column B to a different set of vendors, and move
# the active-cell selection back to I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (between "Alta Equipment Company" and
# "Beacon") to keep the vendor list alphabetically sorted, and fill in
# the new vendor's cells (A: name, G: the "x" column that is filled for
# every vendor row).
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Atlas Crane Service"
$ws.Range("G4").Value = "x"

# Clear out the old "Need Envelope" (column B) marks - they used to sit
# on Garlock Chicago Inc. (now row 8), Metal-Era (row 12),
# Pro Fastening Systems Inc. (row 14) and Stevenson Crane (row 16) after
# the insert shifted everything down by one.
$ws.Range("B8").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B16").ClearContents()

# Put the new "Need Envelope" marks on ABC Supply (row 2), Beacon (row 5)
# and Gemco Supply (row 9).
$ws.Range("B2").Value = "x"
$ws.Range("B5").Value = "x"
$ws.Range("B9").Value = "x"

# Refresh the "last retrieved" date serial in I1.
$ws.Range("I1").Value = 33396

# Move the selection from I2 to I1.
$ws.Range("I1").Select()
